$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2582.6667
$ws.Range("J32").Value = 2937
$ws.Range("L32").Value = 2937
$ws.Range("N32").Value = -3589
$ws.Range("H70").Value = 1356455.6
$ws.Range("J70").Value = 1622.5
$ws.Range("L70").Value = 4867.5
$ws.Range("N70").Value = -5407.5
$ws.Range("H73").Value = 1356455.6
$ws.Range("J73").Value = 1622.5
$ws.Range("L73").Value = 4867.5
$ws.Range("N73").Value = -6739.5
$ws.Range("H98").Value = 1462.6666
$ws.Range("I98").Value = 931.2857
$ws.Range("J98").Value = 3322.5
$ws.Range("K98").Value = 931.2857
$ws.Range("L98").Value = 3322.5
$ws.Range("M98").Value = 566.7143
$ws.Range("N98").Value = -6318.5
$ws.Range("H103").Value = 21740012
$ws.Range("J103").Value = 45455584
$ws.Range("L103").Value = 136366752
$ws.Range("N103").Value = -136367924
$ws.Range("H122").Value = 1462.6666
$ws.Range("I122").Value = 931.2857
$ws.Range("J122").Value = 3322.5
$ws.Range("K122").Value = 2793.8571
$ws.Range("L122").Value = 9967.5
$ws.Range("M122").Value = -343.8571000000002
$ws.Range("N122").Value = -14867.5
$ws.Range("H137").Value = 1223.591
$ws.Range("I137").Value = 1217.6666
$ws.Range("J137").Value = 1348
$ws.Range("K137").Value = 3652.9998
$ws.Range("L137").Value = 4044
$ws.Range("M137").Value = -1102.9998
$ws.Range("N137").Value = -9144
$ws.Range("H138").Value = 4841.2856
$ws.Range("I138").Value = 2981.2632
$ws.Range("J138").Value = 7050.0625
$ws.Range("K138").Value = 8943.7896
$ws.Range("L138").Value = 21150.1875
$ws.Range("M138").Value = -3803.7896
$ws.Range("N138").Value = -31430.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4277.0713
$ws.Range("I2").Value = 4242
$ws.Range("J2").Value = 4487.5
$ws.Range("K2").Value = 4242
$ws.Range("L2").Value = 4487.5
$ws.Range("M2").Value = -4129
$ws.Range("N2").Value = -4713.5
$ws.Range("H32").Value = 2087.8525
$ws.Range("I32").Value = 2039.983
$ws.Range("K32").Value = 2039.983
$ws.Range("M32").Value = -1752.983
$ws.Range("H45").Value = 5347.5
$ws.Range("I45").Value = 4890
$ws.Range("K45").Value = 4890
$ws.Range("M45").Value = -4513
$ws.Range("H61").Value = 1056555.1
$ws.Range("I61").Value = 2677.5386
$ws.Range("J61").Value = 3339956.8
$ws.Range("K61").Value = 2677.5386
$ws.Range("L61").Value = 3339956.8
$ws.Range("M61").Value = -2465.5386
$ws.Range("N61").Value = -3340380.8
$ws.Range("H102").Value = 2107.0588
$ws.Range("J102").Value = 2918.6667
$ws.Range("L102").Value = 2918.6667
$ws.Range("N102").Value = -6162.6667
$ws.Range("H116").Value = 4277.0713
$ws.Range("I116").Value = 4242
$ws.Range("J116").Value = 4487.5
$ws.Range("K116").Value = 4242
$ws.Range("L116").Value = 4487.5
$ws.Range("M116").Value = -1948
$ws.Range("N116").Value = -9075.5
$ws.Range("H132").Value = 3574588.2
$ws.Range("I132").Value = 2746.4707
$ws.Range("J132").Value = 9094708
$ws.Range("K132").Value = 8239.4121
$ws.Range("L132").Value = 27284124
$ws.Range("M132").Value = -5709.4121
$ws.Range("N132").Value = -27289184
$ws.Range("H136").Value = 1056555.1
$ws.Range("I136").Value = 2677.5386
$ws.Range("J136").Value = 3339956.8
$ws.Range("K136").Value = 8032.6158
$ws.Range("L136").Value = 10019870.4
$ws.Range("M136").Value = -5482.6158
$ws.Range("N136").Value = -10024970.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4277.0713
$ws.Range("I3").Value = 4242
$ws.Range("J3").Value = 4487.5
$ws.Range("K3").Value = 4242
$ws.Range("L3").Value = 4487.5
$ws.Range("M3").Value = -4128
$ws.Range("N3").Value = -4715.5
$ws.Range("H86").Value = 4086.7646
$ws.Range("I86").Value = 1429.5
$ws.Range("J86").Value = 5536.1816
$ws.Range("K86").Value = 1429.5
$ws.Range("L86").Value = 5536.1816
$ws.Range("M86").Value = -306.5
$ws.Range("N86").Value = -7782.1816
$ws.Range("H89").Value = 4086.7646
$ws.Range("I89").Value = 1429.5
$ws.Range("J89").Value = 5536.1816
$ws.Range("K89").Value = 7147.5
$ws.Range("L89").Value = 27680.908
$ws.Range("M89").Value = -1531.5
$ws.Range("N89").Value = -38912.908
$ws.Range("H99").Value = 1971.75
$ws.Range("I99").Value = 1768.3
$ws.Range("K99").Value = 1768.3
$ws.Range("M99").Value = -270.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1200.4286
$ws.Range("I22").Value = 675.25
$ws.Range("K22").Value = 675.25
$ws.Range("M22").Value = -325.25
$ws.Range("H31").Value = 3066.5293
$ws.Range("I31").Value = 3242.6155
$ws.Range("J31").Value = 2494.25
$ws.Range("K31").Value = 3242.6155
$ws.Range("L31").Value = 2494.25
$ws.Range("M31").Value = -2947.6155
$ws.Range("N31").Value = -3084.25
$ws.Range("H34").Value = 3066.5293
$ws.Range("I34").Value = 3242.6155
$ws.Range("J34").Value = 2494.25
$ws.Range("K34").Value = 3242.6155
$ws.Range("L34").Value = 2494.25
$ws.Range("M34").Value = -3040.6155
$ws.Range("N34").Value = -2898.25
$ws.Range("H58").Value = 2083.375
$ws.Range("I58").Value = 1374.7778
$ws.Range("J58").Value = 2994.4285
$ws.Range("K58").Value = 1374.7778
$ws.Range("L58").Value = 2994.4285
$ws.Range("M58").Value = -1171.7778
$ws.Range("N58").Value = -3400.4285
$ws.Range("H60").Value = 24768.1
$ws.Range("I60").Value = 9474.4
$ws.Range("J60").Value = 40061.8
$ws.Range("K60").Value = 9474.4
$ws.Range("L60").Value = 40061.8
$ws.Range("M60").Value = -8963.4
$ws.Range("N60").Value = -41083.8
$ws.Range("H86").Value = 8840.736999999999
$ws.Range("I86").Value = 10258.77
$ws.Range("K86").Value = 10258.77
$ws.Range("M86").Value = -9135.77
$ws.Range("H89").Value = 8840.736999999999
$ws.Range("I89").Value = 10258.77
$ws.Range("K89").Value = 51293.85000000001
$ws.Range("M89").Value = -45677.85000000001
$ws.Range("H105").Value = 2049.8
$ws.Range("J105").Value = 2500
$ws.Range("L105").Value = 2500
$ws.Range("N105").Value = -5994
$ws.Range("H136").Value = 2083.375
$ws.Range("I136").Value = 1374.7778
$ws.Range("J136").Value = 2994.4285
$ws.Range("K136").Value = 4124.3334
$ws.Range("L136").Value = 8983.2855
$ws.Range("M136").Value = -1574.3334
$ws.Range("N136").Value = -14083.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1021.8889
$ws.Range("I8").Value = 1021.8889
$ws.Range("K8").Value = 3065.6667
$ws.Range("M8").Value = -2926.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3788.4167
$ws.Range("I102").Value = 4496.5
$ws.Range("K102").Value = 4496.5
$ws.Range("M102").Value = -2874.5
$ws.Range("H126").Value = 2326.75
$ws.Range("I126").Value = 1853.2222
$ws.Range("J126").Value = 3747.3333
$ws.Range("K126").Value = 5559.6666
$ws.Range("L126").Value = 11241.9999
$ws.Range("M126").Value = -3089.6666
$ws.Range("N126").Value = -16181.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5211103
$ws.Range("I68").Value = 8335388.5
$ws.Range("J68").Value = 3961
$ws.Range("K68").Value = 8335388.5
$ws.Range("L68").Value = 3961
$ws.Range("M68").Value = -8334639.5
$ws.Range("N68").Value = -5459
$ws.Range("H71").Value = 5211103
$ws.Range("I71").Value = 8335388.5
$ws.Range("J71").Value = 3961
$ws.Range("K71").Value = 41676942.5
$ws.Range("L71").Value = 19805
$ws.Range("M71").Value = -41673198.5
$ws.Range("N71").Value = -27293
$ws.Range("H93").Value = 6179285
$ws.Range("I93").Value = 3334.3333
$ws.Range("J93").Value = 9267260
$ws.Range("K93").Value = 3334.3333
$ws.Range("L93").Value = 9267260
$ws.Range("M93").Value = -2086.3333
$ws.Range("N93").Value = -9269756

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3387.3428
$ws.Range("I107").Value = 1593.95
$ws.Range("J107").Value = 5778.533
$ws.Range("K107").Value = 4781.85
$ws.Range("L107").Value = 17335.599
$ws.Range("M107").Value = -2861.85
$ws.Range("N107").Value = -21175.599
$ws.Range("H113").Value = 946.2308
$ws.Range("I113").Value = 900.15
$ws.Range("J113").Value = 1099.8334
$ws.Range("K113").Value = 2700.45
$ws.Range("L113").Value = 3299.5002
$ws.Range("M113").Value = -530.4499999999998
$ws.Range("N113").Value = -7639.5002
$ws.Range("H132").Value = 386548.12
$ws.Range("I132").Value = 1968.3
$ws.Range("K132").Value = 5904.9
$ws.Range("M132").Value = -3374.9
